$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1785.7142
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 1780
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 1780
$ws.Range("M40").Value = -1625
$ws.Range("N40").Value = -2130

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9300.478999999999
$ws.Range("I76").Value = 12869.77
$ws.Range("J76").Value = 4660.4
$ws.Range("K76").Value = 12869.77
$ws.Range("L76").Value = 4660.4
$ws.Range("M76").Value = -12554.77
$ws.Range("N76").Value = -5290.4

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 9300.478999999999
$ws.Range("I79").Value = 12869.77
$ws.Range("J79").Value = 4660.4
$ws.Range("K79").Value = 12869.77
$ws.Range("L79").Value = 4660.4
$ws.Range("M79").Value = -11777.77
$ws.Range("N79").Value = -6844.4

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2489.5
$ws.Range("I106").Value = 2124.375
$ws.Range("J106").Value = 3950
$ws.Range("K106").Value = 2124.375
$ws.Range("L106").Value = 3950
$ws.Range("M106").Value = -1493.375
$ws.Range("N106").Value = -5212

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 83336750
$ws.Range("I113").Value = 166668340
$ws.Range("J113").Value = 5162.6665
$ws.Range("K113").Value = 166668340
$ws.Range("L113").Value = 5162.6665
$ws.Range("M113").Value = -166665086
$ws.Range("N113").Value = -11670.6665

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3395.842
$ws.Range("I2").Value = 1310.091
$ws.Range("J2").Value = 6263.75
$ws.Range("K2").Value = 1310.091
$ws.Range("L2").Value = 6263.75
$ws.Range("M2").Value = -1197.091
$ws.Range("N2").Value = -6489.75

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5500
$ws.Range("I45").Value = 6103
$ws.Range("J45").Value = 5098
$ws.Range("K45").Value = 6103
$ws.Range("L45").Value = 5098
$ws.Range("M45").Value = -5726
$ws.Range("N45").Value = -5852

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2498.75
$ws.Range("I63").Value = 2053.3333
$ws.Range("J63").Value = 3071.4285
$ws.Range("K63").Value = 2053.3333
$ws.Range("L63").Value = 3071.4285
$ws.Range("M63").Value = -1367.3333
$ws.Range("N63").Value = -4443.4285

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2498.75
$ws.Range("I66").Value = 2053.3333
$ws.Range("J66").Value = 3071.4285
$ws.Range("K66").Value = 10266.6665
$ws.Range("L66").Value = 15357.1425
$ws.Range("M66").Value = -6834.666499999999
$ws.Range("N66").Value = -22221.1425

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2962.0857
$ws.Range("I97").Value = 2000.6111
$ws.Range("J97").Value = 3980.1177
$ws.Range("K97").Value = 2000.6111
$ws.Range("L97").Value = 3980.1177
$ws.Range("M97").Value = -1504.6111
$ws.Range("N97").Value = -4972.1177

# ARM row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 52947.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 52947.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 52947.5
$ws.Range("N113").Value = -61625.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3395.842
$ws.Range("I116").Value = 1310.091
$ws.Range("J116").Value = 6263.75
$ws.Range("K116").Value = 1310.091
$ws.Range("L116").Value = 6263.75
$ws.Range("M116").Value = 983.9090000000001
$ws.Range("N116").Value = -10851.75

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3395.842
$ws.Range("I3").Value = 1310.091
$ws.Range("J3").Value = 6263.75
$ws.Range("K3").Value = 1310.091
$ws.Range("L3").Value = 6263.75
$ws.Range("M3").Value = -1196.091
$ws.Range("N3").Value = -6491.75

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3559.6924
$ws.Range("I107").Value = 2388.818
$ws.Range("J107").Value = 9999.5
$ws.Range("K107").Value = 2388.818
$ws.Range("L107").Value = 9999.5
$ws.Range("M107").Value = -468.8180000000002
$ws.Range("N107").Value = -13839.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1914.0303
$ws.Range("I134").Value = 1936.6552
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 5809.9656
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -3274.9656
$ws.Range("N134").Value = -10320

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2336.6765
$ws.Range("I31").Value = 1913.1538
$ws.Range("J31").Value = 2598.8572
$ws.Range("K31").Value = 1913.1538
$ws.Range("L31").Value = 2598.8572
$ws.Range("M31").Value = -1618.1538
$ws.Range("N31").Value = -3188.8572

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2336.6765
$ws.Range("I34").Value = 1913.1538
$ws.Range("J34").Value = 2598.8572
$ws.Range("K34").Value = 1913.1538
$ws.Range("L34").Value = 2598.8572
$ws.Range("M34").Value = -1711.1538
$ws.Range("N34").Value = -3002.8572

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1519.9
$ws.Range("I99").Value = 1314.1428
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1314.1428
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 183.8571999999999
$ws.Range("N99").Value = -4996

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1519.9
$ws.Range("I126").Value = 1314.1428
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3942.4284
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1472.4284
$ws.Range("N126").Value = -10940

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3103290.5
$ws.Range("I24").Value = 10333335
$ws.Range("J24").Value = 4699.857
$ws.Range("K24").Value = 10333335
$ws.Range("L24").Value = 4699.857
$ws.Range("M24").Value = -10333162
$ws.Range("N24").Value = -5045.857

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3563.3333
$ws.Range("I80").Value = 3160
$ws.Range("J80").Value = 4101.1113
$ws.Range("K80").Value = 3160
$ws.Range("L80").Value = 4101.1113
$ws.Range("M80").Value = -2162
$ws.Range("N80").Value = -6097.1113

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3563.3333
$ws.Range("I83").Value = 3160
$ws.Range("J83").Value = 4101.1113
$ws.Range("K83").Value = 15800
$ws.Range("L83").Value = 20505.5565
$ws.Range("M83").Value = -10808
$ws.Range("N83").Value = -30489.5565

# GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3400
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -4272

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 1000
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("M48").Value = -339

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3724.1667
$ws.Range("I61").Value = 2648.75
$ws.Range("J61").Value = 5875
$ws.Range("K61").Value = 2648.75
$ws.Range("L61").Value = 5875
$ws.Range("M61").Value = -2446.75
$ws.Range("N61").Value = -6279

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3724.1667
$ws.Range("I113").Value = 2648.75
$ws.Range("J113").Value = 5875
$ws.Range("K113").Value = 2648.75
$ws.Range("L113").Value = 5875
$ws.Range("M113").Value = -478.75
$ws.Range("N113").Value = -10215

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7542.0586
$ws.Range("I81").Value = 17718.834
$ws.Range("J81").Value = 1991.091
$ws.Range("K81").Value = 35437.668
$ws.Range("L81").Value = 3982.182
$ws.Range("M81").Value = -34376.668
$ws.Range("N81").Value = -6104.182

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7542.0586
$ws.Range("I84").Value = 17718.834
$ws.Range("J84").Value = 1991.091
$ws.Range("K84").Value = 177188.34
$ws.Range("L84").Value = 19910.91
$ws.Range("M84").Value = -171884.34
$ws.Range("N84").Value = -30518.91

Write-Host "Edit complete"
